# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" stats sheet: bump the "last updated" timestamp,
# swap two pairs of countries whose ranking order changed (Costa Rica now
# ahead of Venezuela; Timor Oriental now ahead of Santa Lucia), and update the
# numeric columns (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for every country whose counts moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 23:02"

# --- Re-ranked country pairs (labels swap, row 54/55 also gets new data) -
$ws.Range("A54").Value = "Costa Rica"
$ws.Range("A55").Value = "Venezuela"

$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Updated numeric data ----------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 7347285
$ws.Range("C4").Value = 25942
$ws.Range("D4").Value = 4597165
$ws.Range("E4").Value = 2540441
$ws.Range("G4").Value = 226
$ws.Range("H4").Value = 209679

# Row 5: India
$ws.Range("B5").Value = 6143019
$ws.Range("C5").Value = 69671
$ws.Range("D5").Value = 5098573
$ws.Range("E5").Value = 948095
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 96351

# Row 13: Sudafrica
$ws.Range("B13").Value = 671669
$ws.Range("C13").Value = 903
$ws.Range("D13").Value = 604478
$ws.Range("E13").Value = 50605
$ws.Range("G13").Value = 188
$ws.Range("H13").Value = 16586

# Row 25: Alemania
$ws.Range("B25").Value = 288617
$ws.Range("C25").Value = 2279
$ws.Range("E25").Value = 28272

# Row 27: Israel
$ws.Range("B27").Value = 233265
$ws.Range("C27").Value = 2239
$ws.Range("D27").Value = 165191
$ws.Range("E27").Value = 66567
$ws.Range("G27").Value = 41
$ws.Range("H27").Value = 1507

# Row 29: Canada
$ws.Range("B29").Value = 154628
$ws.Range("C29").Value = 1503
$ws.Range("D29").Value = 131947
$ws.Range("E29").Value = 13410
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 9271

# Row 30: Ecuador
$ws.Range("B30").Value = 134965
$ws.Range("C30").Value = 218
$ws.Range("E30").Value = 11389
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 11280

# Row 54: now Costa Rica (new data)
$ws.Range("B54").Value = 73714
$ws.Range("C54").Value = 935
$ws.Range("D54").Value = 29420
$ws.Range("E54").Value = 43433
$ws.Range("G54").Value = 19
$ws.Range("H54").Value = 861

# Row 55: now Venezuela (previous Venezuela data, one column unchanged)
$ws.Range("B55").Value = 72691
$ws.Range("D55").Value = 62427
$ws.Range("E55").Value = 9658
$ws.Range("H55").Value = 606

# Row 56: Barein
$ws.Range("E56").Value = 6229
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 245

# Row 85: Costa de Marfil
$ws.Range("B85").Value = 19641
$ws.Range("C85").Value = 12
$ws.Range("D85").Value = 19202
$ws.Range("E85").Value = 319

# Row 100: Guinea
$ws.Range("B100").Value = 10598
$ws.Range("C100").Value = 18
$ws.Range("D100").Value = 9940
$ws.Range("E100").Value = 592

# Row 116: Cabo Verde
$ws.Range("B116").Value = 5817
$ws.Range("C116").Value = 46
$ws.Range("E116").Value = 624

# Row 134: Siria
$ws.Range("B134").Value = 4102
$ws.Range("C134").Value = 30
$ws.Range("D134").Value = 1074
$ws.Range("E134").Value = 2834
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 194
